# "upload iris data 3"
# Resize/reposition the two pictures on slide 1 (shrink them to make
# room for a third picture that will be added alongside them later),
# and slide the "(B)" caption textbox left so it stays above the
# (now narrower/relocated) second picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Picture 3 (left iris picture, rId2) ---------------------------------
# before: off (961858,1118937) ext (5944268,3344815)
# after:  off (961858,1118938) ext (3164541,1780674)
$picLeft = $s.Shapes.Item("Picture 3")
$picLeft.Left   = 75.73685455322266
$picLeft.Top    = 88.10536193847656
$picLeft.Width  = 249.1764678955078
$picLeft.Height = 140.21055603027344

# --- Picture 4 (right iris picture, rId3) --------------------------------
# before: off (7319877,1118937) ext (3100773,3344815)
# after:  off (4528316,1118937) ext (1650754,1780674)
$picRight = $s.Shapes.Item("Picture 4")
$picRight.Left   = 356.5603332519531
$picRight.Top    = 88.10527801513672
$picRight.Width  = 129.9806365966797
$picRight.Height = 140.21055603027344

# --- TextBox 6 ("(B)" caption) --------------------------------------------
# before: off (6829925,657272)
# after:  off (4078271,664839)
$capRight = $s.Shapes.Item("TextBox 6")
$capRight.Left = 321.12371826171875
$capRight.Top  = 52.34952926635742

# --- Presentation-level empty slide-guide extension ----------------------
# The authored OOXML also gains an (empty) p15:sldGuideLst extension on
# <p:presentation>. PowerPoint writes this automatically once the Guides
# UI/collection has been touched; reflect that intent here. If the
# Guides object model isn't backed by this host, this is a harmless no-op.
try {
    $null = $p.Guides
} catch {
}
